# "taking it to the top" — apply the edits described by the commit.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title: "...Highest Point In Each..." -> "...Highest Point in Each..."
#    (also removes the now-unwanted w:proofErr gramStart/gramEnd markers,
#    since the whole matched range gets replaced and re-run).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Highest Point In Each", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Highest Point in Each", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Insert the new quote paragraph right after the title paragraph.
#    The paragraph mark formatting of paragraph 1 (centered, bold, 28pt)
#    is inherited automatically by InsertParagraphAfter.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$quotePara = $d.Paragraphs.Item(2)
$quotePara.Range.Text = [char]0x201C + "Sa, I" + [char]0x2019 + "m taking it to the top" + [char]0x201D

# ---------------------------------------------------------------------
# 3. "...visiting the highest point in each of the 50 US States is just"
#    -> "... 50 US states is just"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("the highest point in each of the 50 US States is just", $true, $false, $false, $false, $false, `
    $true, 1, $false, "the highest point in each of the 50 US states is just", 2) | Out-Null

# ---------------------------------------------------------------------
# 4. "...have reached the highest point of each of the 50 US States. From..."
#    -> "...of each of the 50 states. From..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("reached the highest point of each of the 50 US States. From", $true, $false, $false, $false, $false, `
    $true, 1, $false, "reached the highest point of each of the 50 states. From", 2) | Out-Null

# ---------------------------------------------------------------------
# 5. "...another 10 as hills and remaining 7 as landmarks."
#    -> "...hills and the remaining 7 as landmarks."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("hills and remaining 7 as landmarks", $true, $false, $false, $false, $false, `
    $true, 1, $false, "hills and the remaining 7 as landmarks", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. "The views from the parkway were impressive" -> "...are impressive"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("The views from the parkway were impressive", $true, $false, $false, $false, $false, `
    $true, 1, $false, "The views from the parkway are impressive", 2) | Out-Null

# ---------------------------------------------------------------------
# 7. "...happened upon the first few US State highpoints..."
#    -> "...US state highpoints..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("happened upon the first few US State highpoints", $true, $false, $false, $false, $false, `
    $true, 1, $false, "happened upon the first few US state highpoints", 2) | Out-Null

# ---------------------------------------------------------------------
# 8. "...back side of the nearly 900 tall cliffs..." -> "...900 foot tall cliffs..."
#    and relocate the _GoBack bookmark to sit right after "900 foot".
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$d.Content.Find.Execute("nearly 900 tall cliffs", $true, $false, $false, $false, $false, `
    $true, 1, $false, "nearly 900 foot tall cliffs", 2) | Out-Null

$markRange = $d.Content
$markRange.Find.Execute("nearly 900 foot", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$markRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $markRange) | Out-Null
